$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("lat_long")

$ws.Rows.Item(136).Insert()
$ws.Range("B136:D136").Style = "Normal"
$ws.Cells.Item(136, 1).Value = "TNS"
$ws.Cells.Item(136, 2).Value = 722
$ws.Cells.Item(136, 3).Value = 38.255028000000003
$ws.Cells.Item(136, 4).Value = -121.68858299999999

$ws.Activate()
$excel.ActiveWindow.ScrollRow = 100
$ws.Range("C136:D136").Select()
